# Updates the cryptocurrency price/volume table to reflect the latest
# data pull (GitHub Actions scheduled refresh).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Most "Price" values in column D look like plain numbers (e.g. "215.10"),
# but the sheet stores them as TEXT (they use "." as a thousands separator
# in some rows, e.g. "27.070.23", which is not a valid number). To make sure
# Excel keeps every updated Price cell as text - instead of silently
# re-interpreting a cell like "215.11" as the number 215.11 - values that
# would otherwise parse as a plain number are entered with a leading
# apostrophe, Excel's standard "treat this as text" input prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.047.76"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.655.30"
$ws.Range("E3").Value = "  +3.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'215.11"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.12%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +1.76%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.66%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'19.68"
$ws.Range("E10").Value = "  +3.67%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  +1.12%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "1.890.59"
$ws.Range("E12").Value = "  +3.91%  "

# Row 13 - Wrapped Ether
$ws.Range("D13").Value = "1.657.82"
$ws.Range("E13").Value = "  +3.96%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.17%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +2.95%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.22"
$ws.Range("E16").Value = "  +2.75%  "

# Row 17/18 - WrappedBTC and BitcoinCash swapped rank order
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'238.83"
$ws.Range("E17").Value = "  +3.99%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "27.048.10"
$ws.Range("E18").Value = "  +3.01%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'7.83"
$ws.Range("E19").Value = "  +2.27%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +1.26%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.04%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  +4.45%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +4.34%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +3.58%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'145.78"
$ws.Range("E25").Value = "  -0.27%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'7.12"
$ws.Range("E27").Value = "  +2.14%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +1.45%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +3.31%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.69%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.85%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.29"
$ws.Range("E32").Value = "  +3.20%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.515.34"
$ws.Range("E33").Value = "  +2.81%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "'3.06"
$ws.Range("E34").Value = "  +4.29%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +9.41%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.07%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +1.83%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "'0.886"
$ws.Range("E38").Value = "  +8.39%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.67%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'5.96"
$ws.Range("E40").Value = "  +3.58%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.05%  "

# Row 42 - MXToken
$ws.Range("D42").Value = "'2.25"
$ws.Range("E42").Value = "  +4.11%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'65.95"
$ws.Range("E43").Value = "  +9.27%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.797.98"
$ws.Range("E44").Value = "  +3.75%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "'0.776"
$ws.Range("E45").Value = "  +2.81%  "

# Row 46 - WEMIXToken
$ws.Range("D46").Value = "'0.921"
$ws.Range("E46").Value = "  -1.24%  "

# Row 47 - Quant
$ws.Range("E47").Value = "  +1.49%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  +3.16%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "'0.0507"
$ws.Range("E50").Value = "  +1.02%  "

# Row 51 - Algorand
$ws.Range("D51").Value = "'0.0975"
$ws.Range("E51").Value = "  +2.53%  "
